# Apply the abstract-body edit described in the commit diff.
#
# The document currently has 5 paragraphs: a centered title, a centered
# empty paragraph, and three plain empty paragraphs. The diff replaces the
# last (5th, entirely empty) paragraph with one that holds the abstract
# body text, built up out of several runs.

$d = $word.ActiveDocument

# The target (5th) paragraph is the last, still-empty paragraph in the body.
$target = $d.Paragraphs.Item($d.Paragraphs.Count)

# Collapse to the very start of that paragraph, then type the runs in
# sequence so they land one after another, in document order.
$ip = $target.Range
$ip.Collapse(1)   # wdCollapseStart

$runs = @(
    "This project entails the creation of a mobile application that allows users to democratically elect tracks to be ",
    "pushed to a queue",
    " on a host device. Alongside this, a hybrid recommender system consisting of genre and artist inference models is developed to recommend tracks to groups of users",
    ", which can then be voted for and consequently elected.",
    " ",
    "People in a social setting must typically rely on a single device to control music playback, ",
    "that is to say it is unifocal, ",
    "leading to potential conflicts of preference",
    ". ",
    "The application developed [Insert findings from evaluation]"
)

foreach ($chunk in $runs) {
    $ip.InsertAfter($chunk)
    $ip.Collapse(0)   # wdCollapseEnd - move insertion point past what was just typed
}

Write-Host "Final paragraph text:" $target.Range.Text
